$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data table for rows 2-10 (A1-style header stays on row 1)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics

$data = @(
    @("ECs","Vegfc","Vipr2","ECs",3,1,4.425865000000001,13.277595,0.4619841037548157,0.4696223785602887,2,0.6666666666666666,0.3508273333333333,1.052482,0.1248367457344541,0.1598401480061556,1.552714415643334,13.97442974079,0.0576725920937996,0.07506451049607937),
    @("ECs","Vegfc","Vipr2","FAPs",3,1,4.425865000000001,13.277595,0.4619841037548157,0.4696223785602887,3,1,0.6131856666666667,1.839557,0.2181930992387853,0.2793730088930353,2.713876980601667,24.424892825415,0.1008017433973158,0.1311998169418919),
    @("ECs","Vegfc","Vipr2","sCs",3,1,4.425865000000001,13.277595,0.4619841037548157,0.4696223785602887,2,1,1.846276,3.692552,0.6569701550267606,0.5607868431008092,8.171368328740002,49.02820997244001,0.3035097682637004,0.2633580511223174),
    @("FAPs","Vegfc","Vipr2","ECs",3,1,4.686805000000001,14.060415,0.4892217470254038,0.4973103589802793,2,0.6666666666666666,0.3508273333333333,1.052482,0.1248367457344541,0.1598401480061556,1.644259300003333,14.79833370003,0.06107285084117577,0.07949016138440221),
    @("FAPs","Vegfc","Vipr2","FAPs",3,1,4.686805000000001,14.060415,0.4892217470254038,0.4973103589802793,3,1,0.6131856666666667,1.839557,0.2181930992387853,0.2793730088930353,2.873881648461667,25.864934836155,0.1067448091984858,0.1389350913419961),
    @("FAPs","Vegfc","Vipr2","sCs",3,1,4.686805000000001,14.060415,0.4892217470254038,0.4973103589802793,2,1,1.846276,3.692552,0.6569701550267606,0.5607868431008092,8.653135588180001,51.91881352908,0.3214040869857422,0.278885106253881),
    @("sCs","Vegfc","Vipr2","ECs",2,1,0.467454,0.9349080000000001,0.04879414921978045,0.03306726245943202,2,0.6666666666666666,0.3508273333333333,1.052482,0.1248367457344541,0.1598401480061556,0.163995640276,0.983973841656,0.006091302799478745,0.005285476125674007),
    @("sCs","Vegfc","Vipr2","FAPs",2,1,0.467454,0.9349080000000001,0.04879414921978045,0.03306726245943202,3,1,0.6131856666666667,1.839557,0.2181930992387853,0.2793730088930353,0.286636092626,1.719816555756,0.01064654664298365,0.009238100609147233),
    @("sCs","Vegfc","Vipr2","sCs",2,1,0.467454,0.9349080000000001,0.04879414921978045,0.03306726245943202,2,1,1.846276,3.692552,0.6569701550267606,0.5607868431008092,0.863049101304,3.452196405216,0.03205629977731805,0.01854368572461078)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
